$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 23; $row++) {
    $ws.Range("C$row").Value = 45204
}
